# Revert "Powerpoint writer: consolidate text run nodes."
#
# A number of paragraphs in this deck (the slide titles, "Slide N
# (Layout)", and the "an image" / "An image" picture captions) still
# have runs where a trailing space is glued onto the preceding word,
# e.g.
#   <a:r><a:t>Slide </a:t></a:r><a:r><a:t>1 </a:t></a:r><a:r><a:t>(Content)</a:t></a:r>
# The fix splits every "word + trailing space" run into two runs - the
# bare word, and a standalone single-space run - leaving every other run
# (and every <a:rPr/>) untouched, e.g.
#   <a:r><a:t>Slide</a:t></a:r><a:r><a:t> </a:t></a:r>
#   <a:r><a:t>1</a:t></a:r><a:r><a:t> </a:t></a:r>
#   <a:r><a:t>(Content)</a:t></a:r>
#
# The object model here doesn't expose a direct "split this run" call
# (TextRange2.Runs/.Words/.Paragraphs items come back text-less), but
# re-assigning a Characters() sub-range's .Text to itself (a textual
# no-op) makes the host re-flow the paragraph and materialize a fresh
# run boundary exactly at that sub-range's edges, with an empty
# <a:rPr/> on both sides - matching the target XML. Doing this for the
# range that covers just the "word" part (i.e. excluding the trailing
# space) yields precisely the "word" / " " split wanted here.
#
# Because re-assigning .Text on a sub-range re-normalizes everything it
# spans (collapsing back together any splits already made strictly
# inside it), the split points within a paragraph must be applied
# back-to-front (rightmost word first) so earlier splits are never
# swallowed by a later, wider one.
#
# Only paragraphs matching the affected text (slide-number titles and
# the image captions) are touched; other multi-word runs elsewhere in
# the deck (e.g. "Even with some text first, these should:") are left
# exactly as they are, since the diff does not touch them.

function Split-TrailingSpaceRuns($tr) {
    $full = $tr.Text
    $len = $full.Length

    $starts = New-Object System.Collections.ArrayList
    $lens = New-Object System.Collections.ArrayList

    $i = 1
    while ($i -le $len) {
        $ch = $full.Substring($i - 1, 1)
        if ($ch -ne " ") {
            $wordStart = $i
            while (($i -le $len) -and ($full.Substring($i - 1, 1) -ne " ")) {
                $i = $i + 1
            }
            $wordLen = $i - $wordStart
            # Only split off the word if it is immediately followed by a
            # single space character (i.e. it currently ends a run that
            # has a trailing space glued onto it).
            if (($i -le $len) -and ($full.Substring($i - 1, 1) -eq " ")) {
                [void]$starts.Add($wordStart)
                [void]$lens.Add($wordLen)
            }
        } else {
            $i = $i + 1
        }
    }

    for ($k = $starts.Count - 1; $k -ge 0; $k--) {
        $c = $tr.Characters($starts[$k], $lens[$k])
        $c.Text = $c.Text
    }
}

function Test-NeedsSplit($text) {
    if ($text -match "^Slide \d") { return $true }
    if ($text -eq "an image") { return $true }
    if ($text -eq "An image") { return $true }
    return $false
}

$p = $ppt.ActivePresentation

for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $slide = $p.Slides.Item($idx)
    for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
        $shape = $slide.Shapes.Item($si)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if (Test-NeedsSplit $para.Text) {
                Split-TrailingSpaceRuns $para
            }
        }
    }
}
